# Trade #9 closed at 2026-02-16 22:58:10 - base_strategy DOWN +0.000%
# Append a new trade row (row 10) to both the "All Trades" and
# "base_strategy" worksheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "base_strategy")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $row = 10

    $ws.Cells.Item($row, 1).Value = 9

    # "2026-02-16" looks like a date to Excel's auto-detection, so force
    # it to be stored as literal text (matches how the existing Date
    # column cells in this sheet are stored).
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-16"

    $ws.Cells.Item($row, 3).Value = "22:58:10"
    $ws.Cells.Item($row, 4).Value = "base_strategy"
    $ws.Cells.Item($row, 5).Value = "DOWN"
    $ws.Cells.Item($row, 6).Value = 0.5

    # Exit Price / Exit Reason are blank (open trade) but still present
    # as empty text cells, just like the other OPEN trades above them.
    # A leading apostrophe forces an explicit (empty) text entry instead
    # of clearing the cell.
    $ws.Cells.Item($row, 7).Value = "'"

    $ws.Cells.Item($row, 8).Value = "OPEN"
    $ws.Cells.Item($row, 9).Value = 0
    $ws.Cells.Item($row, 10).Value = 0
    $ws.Cells.Item($row, 11).Value = 100
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($row, 16).Value = "'"
    $ws.Cells.Item($row, 17).Value = 0
}
